$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.686.43"
$ws.Range("E2").Value = "  -4.40%  "
$ws.Range("D3").Value = "2.941.86"
$ws.Range("E3").Value = "  -2.36%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "549.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.16%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "130.97"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.80%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.513"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.71%  "
$ws.Range("D9").Value = "2.935.37"
$ws.Range("E9").Value = "  -2.41%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.127"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.21%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "4.77"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.70%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.447"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.85%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000222"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.12%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.95"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.79%  "
$ws.Range("E15").Value = "  +1.24%  "
$ws.Range("D16").Value = "3.425.17"
$ws.Range("E16").Value = "  -2.38%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.90"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +5.71%  "
$ws.Range("D18").Value = "2.934.56"
$ws.Range("E18").Value = "  -2.45%  "
$ws.Range("D19").Value = "57.662.14"
$ws.Range("E19").Value = "  -4.36%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "418.19"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.91%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.27"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.28%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.689"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.50%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.26%  "
$ws.Range("E24").Value = "  +1.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "79.87"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.45%  "
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("E27").Value = "  +0.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.47"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.29%  "
$ws.Range("B29").Value = "ImmutableX"
$ws.Range("C29").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.76%  "
$ws.Range("B30").Value = "RenderToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.45"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.46%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.26"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.60%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.01"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.49%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0970"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.01%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.69"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.93%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.937"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.48%  "
$ws.Range("E36").Value = "  +0.80%  "
$ws.Range("B37").Value = "PEPE"
$ws.Range("C37").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D37").Value = "0.0₃0694"
$ws.Range("E37").Value = "  +3.11%  "
$ws.Range("B38").Value = "OKB"
$ws.Range("C38").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "48.30"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.75%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.75"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.56%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.56"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.88%  "
$ws.Range("E41").Value = "  -0.87%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "377.51"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.63%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0346"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.24%  "
$ws.Range("D44").Value = "2.703.89"
$ws.Range("E44").Value = "  +0.79%  "
$ws.Range("E45").Value = "  +0.06%  "
$ws.Range("E46").Value = "  +1.10%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "122.33"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.48%  "
$ws.Range("E48").Value = "  +1.25%  "
$ws.Range("E49").Value = "  -1.97%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.14"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.73%  "
$ws.Range("E51").Value = "  -0.19%  "
